# Update control flow ppt slides:
#   - Swap the programming-language label / sample code between C# and C++
#     (and Pascal_case .NET-style calls to snake_case ones) across the
#     seven "if/else" slides.
#
# Helper: resolve a shape by name anywhere on a slide - either a
# top-level shape, or nested (at any depth) inside a top-level group.
function Get-ShapeByName {
    param($slide, [string]$name)

    try {
        $direct = $slide.Shapes.Item($name)
        if ($direct) { return $direct }
    } catch {
    }

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Type -eq 6) {
            try {
                $found = $shp.GroupItems.Item($name)
                if ($found) { return $found }
            } catch {
            }
        }
    }

    return $null
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

(Get-ShapeByName $s "TextBox 59").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "C++"
(Get-ShapeByName $s "TextBox 66").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "C++"

$code = (Get-ShapeByName $s "TextBox 67").TextFrame.TextRange
$code.Paragraphs(2,1).Runs(1,1).Text = 'write("What language do you use? ");'
$code.Paragraphs(3,1).Runs(2,1).Text = "read_line"
$code.Paragraphs(4,1).Runs(1,1).Text = 'if (language == "C++")'
$code.Paragraphs(6,1).Runs(1,1).Text = '    write_line("Good choice, C++ is a great language.");'

# ---------------------------------------------------------------------
# Slide 2
# ---------------------------------------------------------------------
$s = $p.Slides.Item(2)

(Get-ShapeByName $s "TextBox 59").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "C++"
(Get-ShapeByName $s "TextBox 66").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "C++"

$code = (Get-ShapeByName $s "TextBox 67").TextFrame.TextRange
$code.Paragraphs(2,1).Runs(1,1).Text = 'write("What language do you use? ");'
$code.Paragraphs(3,1).Runs(2,1).Text = "read_line"
$code.Paragraphs(4,1).Runs(1,1).Text = 'if (language == "C++")'
$code.Paragraphs(6,1).Runs(1,1).Text = '    write_line("Good choice, C++ is a great language.");'

(Get-ShapeByName $s "TextBox 3").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "Good choice, C++ is a great language"

# ---------------------------------------------------------------------
# Slide 3
# ---------------------------------------------------------------------
$s = $p.Slides.Item(3)

(Get-ShapeByName $s "TextBox 94").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "C++"
(Get-ShapeByName $s "TextBox 97").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "C++"

$code = (Get-ShapeByName $s "TextBox 98").TextFrame.TextRange
$code.Paragraphs(1,1).Runs(1,1).Text = 'if (language == "C++")'
$code.Paragraphs(3,1).Runs(1,1).Text = '    write_line("Good choice, C++ is a great language.");'
$code.Paragraphs(7,1).Runs(1,1).Text = '    write_line("Well... good luck with that!");'
$code.Paragraphs(9,1).Runs(1,1).Text = 'write_line("Great chat!");'

(Get-ShapeByName $s "TextBox 3").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "Good choice, C++ is a great language"

# ---------------------------------------------------------------------
# Slide 4
# ---------------------------------------------------------------------
$s = $p.Slides.Item(4)

(Get-ShapeByName $s "TextBox 94").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "C++"
(Get-ShapeByName $s "TextBox 97").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "C++"

$code = (Get-ShapeByName $s "TextBox 98").TextFrame.TextRange
$code.Paragraphs(1,1).Runs(1,1).Text = 'if (language == "C++")'
$code.Paragraphs(3,1).Runs(1,1).Text = '    write_line("Good choice, C++ is a great language.");'
$code.Paragraphs(7,1).Runs(1,1).Text = '    write_line("Well... good luck with that!");'
$code.Paragraphs(9,1).Runs(1,1).Text = 'write_line("Great chat!");'

(Get-ShapeByName $s "TextBox 3").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "Good choice, C++ is a great language"

# ---------------------------------------------------------------------
# Slide 5  (title reverts C++ -> C#, code becomes C++/write_line)
# ---------------------------------------------------------------------
$s = $p.Slides.Item(5)

(Get-ShapeByName $s "TextBox 94").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "C#"
(Get-ShapeByName $s "TextBox 97").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "C#"

$code = (Get-ShapeByName $s "TextBox 98").TextFrame.TextRange
$code.Paragraphs(1,1).Runs(1,1).Text = 'if (language == "C++")'
$code.Paragraphs(3,1).Runs(1,1).Text = '    write_line("Good choice, C++ is a great language.");'
$code.Paragraphs(7,1).Runs(1,1).Text = '    write_line("Well... good luck with that!");'
$code.Paragraphs(9,1).Runs(1,1).Text = 'write_line("Great chat!");'

# ---------------------------------------------------------------------
# Slide 6  (title reverts C++ -> C#, code becomes C++/write_line)
# ---------------------------------------------------------------------
$s = $p.Slides.Item(6)

(Get-ShapeByName $s "TextBox 94").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "C#"
(Get-ShapeByName $s "TextBox 97").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "C#"

$code = (Get-ShapeByName $s "TextBox 98").TextFrame.TextRange
$code.Paragraphs(1,1).Runs(1,1).Text = 'if (language == "C++")'
$code.Paragraphs(3,1).Runs(1,1).Text = '    write_line("Good choice, C++ is a great language.");'
$code.Paragraphs(7,1).Runs(1,1).Text = '    write_line("Well... good luck with that!");'
$code.Paragraphs(9,1).Runs(1,1).Text = 'write_line("Great chat!");'

# ---------------------------------------------------------------------
# Slide 7  (title reverts C++ -> C#, code becomes C++/write_line;
#           paragraph 7 has "WriteLine" split into its own run)
# ---------------------------------------------------------------------
$s = $p.Slides.Item(7)

(Get-ShapeByName $s "TextBox 94").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "C#"
(Get-ShapeByName $s "TextBox 97").TextFrame.TextRange.Paragraphs(1,1).Runs(1,1).Text = "C#"

$code = (Get-ShapeByName $s "TextBox 98").TextFrame.TextRange
$code.Paragraphs(1,1).Runs(1,1).Text = 'if (language == "C++")'
$code.Paragraphs(3,1).Runs(1,1).Text = '    write_line("Good choice, C++ is a great language.");'
$code.Paragraphs(7,1).Runs(1,1).Text = "    write_line"
$code.Paragraphs(9,1).Runs(1,1).Text = 'write_line("Great chat!");'
